# PM12 Tidsregistrering for Emil.xlsx - add newly logged time-tracking rows
# (Krydstjek 04 / Review ATD07b / Review DD07 / Ret ENV08 move-rename-vejledning)
# to the "Ark1" time sheet, rows 22-26 (columns A-F). Columns G (hours) and H
# (running total) already contain formulas that recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# --- Row 22: Krydstjek 04 / Reviewer -----------------------------------
$ws.Range("A22").Value = "Krydstjek 04"
$ws.Range("B22").Value = "Reviewer"
$ws.Range("C22").Value = 43889
$ws.Range("D22").Value = 0.3888888888888889
$ws.Range("E22").Value = 0.42708333333333331
$ws.Range("F22").Value = 0.041666666666666664
$ws.Range("F22").NumberFormat = "h:mm"

# --- Row 23: Review ATD07b / Reviewer -----------------------------------
$ws.Range("A23").Value = "Review ATD07b"
$ws.Range("B23").Value = "Reviewer"
$ws.Range("C23").Value = 43889
$ws.Range("D23").Value = 0.45833333333333331
$ws.Range("E23").Value = 0.47222222222222227
$ws.Range("F23").Value = 0.0069444444444444441
$ws.Range("F23").NumberFormat = "h:mm"

# --- Row 24: Review DD07 / Reviewer -------------------------------------
$ws.Range("A24").Value = "Review DD07"
$ws.Range("B24").Value = "Reviewer"
$ws.Range("C24").Value = 43889
$ws.Range("D24").Value = 0.48958333333333331
$ws.Range("E24").Value = 0.49305555555555558
$ws.Range("F24").Value = 0.0069444444444444441
$ws.Range("F24").NumberFormat = "h:mm"

# --- Row 25: Review ATD07b / Reviewer -----------------------------------
$ws.Range("A25").Value = "Review ATD07b"
$ws.Range("B25").Value = "Reviewer"
$ws.Range("C25").Value = 43889
$ws.Range("D25").Value = 0.52083333333333337
$ws.Range("E25").Value = 0.53472222222222221
$ws.Range("F25").Value = 0.0069444444444444441
$ws.Range("F25").NumberFormat = "h:mm"
$ws.Range("F25").HorizontalAlignment = -4108

# --- Row 26: Ret ENV08 move-rename-vejledning / Tool Specialist ---------
$ws.Range("A26").Value = "Ret ENV08 move-rename-vejledning"
$ws.Range("B26").Value = "Tool Specialist"
$ws.Range("C26").Value = 43889
$ws.Range("D26").Value = 0.53472222222222221
$ws.Range("E26").Value = 0.54166666666666663
$ws.Range("F26").Value = 0.0069444444444444441
$ws.Range("F26").NumberFormat = "h:mm"
$ws.Range("F26").HorizontalAlignment = -4108

# --- Restore the view state as left by the author (best effort): scroll so
# row 6 is the top visible row, and leave the active cell on C27. ---------
$win = $excel.ActiveWindow
$win.ScrollRow = 6
$win.ScrollColumn = 1
$ws.Range("C27").Select()

Write-Host "Rows 22-26 populated; totals recalculated through row 32."
